$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = "Should"
$ws.Range("E2").Value = "play starts/resumes; pause toggles; stop halts & resets; errors handled; state reflected in UI"
$ws.Range("E3").Value = "next advances; prev goes back; wrap behavior; single/empty list handling"
$ws.Range("E4").Value = "shows title/artist/duration; updates on change; fallback to filename; formatting"
$ws.Range("E6").Value = "displays elapsed/total; updates each second; paused stops increment; unknown duration handled"
$ws.Range("E7").Value = "text bar matches position; parse mm:ss or seconds; clamp to range; state preserved"
$ws.Range("E8").Value = "maps keys; works when UI focused; conflicts handled; repeat-safe"
$ws.Range("E9").Value = "ff adds 5s; rw subtracts 5s; clamps; works in paused/playing"
$ws.Range("E10").Value = "toggle mute; remember previous volume; set volume while muted updates stored value"
$ws.Range("E11").Value = "current track prefixed icon; updates on state changes; accurate index mapping"
$ws.Range("E12").Value = "help lists commands & descriptions; help <cmd> shows detail; unknown suggests help"
$ws.Range("E13").Value = "separate playback thread/task; CLI non-blocking; graceful shutdown; resource cleanup"
$ws.Range("E14").Value = "create playlist; rename existing; delete by name/id; errors for duplicates/not found"
$ws.Range("E15").Value = "add by id/path; remove by index/id; reorder by move/swap; persistence"
$ws.Range("E16").Value = "case-insensitive search; partial matches; multiple fields; empty result handled"
$ws.Range("E17").Value = "tabular layout; columns aligned; pagination if needed; null-safe values"
$ws.Range("E18").Value = "list with name & count; select by name/number; errors for invalid"
$ws.Range("E19").Value = "show numbered tracks with durations; handles empty; updates on changes"
$ws.Range("E20").Value = "add by search result/id; confirm existence; prevent duplicates (or allow as option)"
$ws.Range("E21").Value = "success messages include playlist & song; failures include reason; no silent errors"
$ws.Range("E22").Value = "scan path recursively; supported extensions only; update library; handle duplicates"
$ws.Range("E23").Value = "shows track count and hh:mm:ss; recalculates on change; accurate sum"
$ws.Range("E24").Value = "merge A+B; optional remove duplicates; stable order; returns new playlist"
$ws.Range("E25").Value = "deep copy retains order; new name required; independent edits"
$ws.Range("E26").Value = "toggle shuffle; deterministic with seed; persists state; respects queue"
$ws.Range("E27").Value = "loop single track; loop playlist; off state; interactions with next/prev"
$ws.Range("E28").Value = "show current queue order; show history list; empty states handled"
$ws.Range("E29").Value = "enqueue by id; remove by position/id; no duplicates optional; immediate effect"
$ws.Range("E30").Value = "place song at position 1 in queue; respects currently playing; id validation"
$ws.Range("E31").Value = "one command clears queue; playback unaffected until next track; confirm message"
$ws.Range("E34").Value = "list all liked tracks; sort by date liked or title; empty state handled"
$ws.Range("E35").Value = "sort by chosen key; stable sort; reversible; persists order"
$ws.Range("E36").Value = "track play counts; list top N; ties deterministic"
$ws.Range("E37").Value = "set timer minutes; cancels/overrides; stops playback on expiry; persistence optional"
$ws.Range("E38").Value = "persist on exit; load on start; defaults if file missing; error-tolerant"
$ws.Range("E39").Value = "define schedules; next-run calculation; execute play at time; overlap handling"
$ws.Range("E40").Value = "checkpoint periodically; restore on start; file corruption safe"
$ws.Range("E42").Value = "create tags; assign to songs; filter query (AND/OR/basic); persist"
$ws.Range("E43").Value = "list by added timestamp; configurable window; empty state handled"
$ws.Range("E44").Value = "create/switch profiles; isolated settings/library/playlists; default profile"
$ws.Range("E45").Value = "accumulate per track/artist; total listening time; render stats view"
$ws.Range("E47").Value = "export fields index/title/artist/duration/path; write CSV; handle IO errors"
$ws.Range("E5").Value = "accepts numeric input; clamps 0-100; errors for invalid; immediate effect"
$ws.Range("E49").Value = "set rating 1-10; filter by rating; auto playlists per rating"
$ws.Range("E48").Value = "edit fields; validation; data persists"
$ws.Range("E46").Value = "add/remove changed files; robustness to rapid changes"
$ws.Range("E41").Value = "import local/external; skip unsupported; duplicate; progress reporting"
$ws.Range("E33").Value = "toggle like; persist per track; "
$ws.Range("E32").Value = "accept 0.5-2.0; clamp; persists; affects playback timing"
$ws.Range("I2").Value = "Backlog"
$ws.Range("I3").Value = "Backlog"
$ws.Range("I4").Value = "Backlog"
$ws.Range("I5").Value = "Backlog"
$ws.Range("D6").Value = "Must"
$ws.Range("I6").Value = "Backlog"
$ws.Range("D7").Value = "Must"
$ws.Range("I7").Value = "Backlog"
$ws.Range("I8").Value = "Backlog"
$ws.Range("D9").Value = "Should"
$ws.Range("I9").Value = "Backlog"
$ws.Range("D10").Value = "Must"
$ws.Range("I10").Value = "Backlog"
$ws.Range("D11").Value = "Should"
$ws.Range("I11").Value = "Backlog"
$ws.Range("D12").Value = "Must"
$ws.Range("I12").Value = "Backlog"
$ws.Range("D13").Value = "Must"
$ws.Range("I13").Value = "Backlog"
$ws.Range("D14").Value = "Must"
$ws.Range("I14").Value = "Backlog"
$ws.Range("D15").Value = "Must"
$ws.Range("I15").Value = "Backlog"
$ws.Range("D16").Value = "Must"
$ws.Range("I16").Value = "Backlog"
$ws.Range("D17").Value = "Should"
$ws.Range("I17").Value = "Backlog"
$ws.Range("D18").Value = "Must"
$ws.Range("I18").Value = "Backlog"
$ws.Range("D19").Value = "Must"
$ws.Range("I19").Value = "Backlog"
$ws.Range("D20").Value = "Should"
$ws.Range("I20").Value = "Backlog"
$ws.Range("D21").Value = "Should"
$ws.Range("I21").Value = "Backlog"
$ws.Range("D22").Value = "Must"
$ws.Range("I22").Value = "Backlog"
$ws.Range("D23").Value = "Should"
$ws.Range("I23").Value = "Backlog"
$ws.Range("D24").Value = "Should"
$ws.Range("I24").Value = "Backlog"
$ws.Range("D25").Value = "Should"
$ws.Range("I25").Value = "Backlog"
$ws.Range("D26").Value = "Must"
$ws.Range("I26").Value = "Backlog"
$ws.Range("D27").Value = "Must"
$ws.Range("I27").Value = "Backlog"
$ws.Range("D28").Value = "Should"
$ws.Range("I28").Value = "Backlog"
$ws.Range("D29").Value = "Must"
$ws.Range("I29").Value = "Backlog"
$ws.Range("D30").Value = "Should"
$ws.Range("I30").Value = "Backlog"
$ws.Range("D31").Value = "Should"
$ws.Range("I31").Value = "Backlog"
$ws.Range("D32").Value = "Should"
$ws.Range("I32").Value = "Backlog"
$ws.Range("D33").Value = "Should"
$ws.Range("I33").Value = "Backlog"
$ws.Range("D34").Value = "Should"
$ws.Range("I34").Value = "Backlog"
$ws.Range("D35").Value = "Should"
$ws.Range("I35").Value = "Backlog"
$ws.Range("D36").Value = "Should"
$ws.Range("I36").Value = "Backlog"
$ws.Range("D37").Value = "Should"
$ws.Range("I37").Value = "Backlog"
$ws.Range("D38").Value = "Must"
$ws.Range("I38").Value = "Backlog"
$ws.Range("D39").Value = "Should"
$ws.Range("I39").Value = "Backlog"
$ws.Range("D40").Value = "Must"
$ws.Range("I40").Value = "Backlog"
$ws.Range("D41").Value = "Must"
$ws.Range("I41").Value = "Backlog"
$ws.Range("D42").Value = "Should"
$ws.Range("I42").Value = "Backlog"
$ws.Range("D43").Value = "Should"
$ws.Range("I43").Value = "Backlog"
$ws.Range("D44").Value = "Should"
$ws.Range("I44").Value = "Backlog"
$ws.Range("D45").Value = "Should"
$ws.Range("I45").Value = "Backlog"
$ws.Range("D46").Value = "Should"
$ws.Range("I46").Value = "Backlog"
$ws.Range("D47").Value = "Should"
$ws.Range("I47").Value = "Backlog"
$ws.Range("D48").Value = "Should"
$ws.Range("I48").Value = "Backlog"
$ws.Range("D49").Value = "Should"
$ws.Range("I49").Value = "Backlog"
